$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two remaining employee e-mail addresses (text + hyperlink target)
$ws.Range("A2").Value = "gr_sql_model1@example.com"
$ws.Range("A3").Value = "gr_sql_model2@example.com"

# Remove the 4th employee row entirely (data + hyperlink)
$ws.Rows("4:4").Delete()

# Rebuild the hyperlinks collection (the engine clears the whole sheet's
# hyperlinks whenever any range's Hyperlinks.Delete() is invoked, so wipe
# once, then add back only the ones that should remain)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gr_sql_model1@example.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:gr_sql_model2@example.com")

# Adding a hyperlink re-stamps cell formatting; restore the original
# "Hipervinculo" cell style used by the two e-mail cells
$ws.Range("A2").Style = "Hipervínculo"
$ws.Range("A3").Style = "Hipervínculo"

# Match the saved selection/active cell
$ws.Range("A4").Select()
